$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 922.7655999999999
$ws.Range("J17").Value = 922.7655999999999
$ws.Range("L17").Value = 2768.2968
$ws.Range("N17").Value = -3104.2968
# Row 28
$ws.Range("H28").Value = 2720560.2
$ws.Range("J28").Value = 21499.5
$ws.Range("L28").Value = 21499.5
$ws.Range("N28").Value = -22469.5
# Row 125
$ws.Range("H125").Value = 2012
$ws.Range("J125").Value = 2068
$ws.Range("L125").Value = 18612
$ws.Range("N125").Value = -23532
# Row 138
$ws.Range("H138").Value = 2812.82
$ws.Range("J138").Value = 3346.238
$ws.Range("L138").Value = 10038.714
$ws.Range("N138").Value = -20318.714

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2891.012
$ws.Range("I32").Value = 2244.4365
$ws.Range("K32").Value = 2244.4365
$ws.Range("M32").Value = -1957.4365
# Row 61
$ws.Range("H61").Value = 5727
$ws.Range("I61").Value = 3699.6667
$ws.Range("J61").Value = 7754.3335
$ws.Range("K61").Value = 3699.6667
$ws.Range("L61").Value = 7754.3335
$ws.Range("M61").Value = -3487.6667
$ws.Range("N61").Value = -8178.3335
# Row 74
$ws.Range("H74").Value = 1343.2188
$ws.Range("I74").Value = 1107.6818
$ws.Range("J74").Value = 1861.4
$ws.Range("K74").Value = 1107.6818
$ws.Range("L74").Value = 1861.4
$ws.Range("M74").Value = -233.6818000000001
$ws.Range("N74").Value = -3609.4
# Row 77
$ws.Range("H77").Value = 1343.2188
$ws.Range("I77").Value = 1107.6818
$ws.Range("J77").Value = 1861.4
$ws.Range("K77").Value = 5538.409000000001
$ws.Range("L77").Value = 9307
$ws.Range("M77").Value = -1170.409000000001
$ws.Range("N77").Value = -18043
# Row 110
$ws.Range("H110").Value = 1562.2
$ws.Range("I110").Value = 1562.2
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1562.2
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 482.8
$ws.Range("N110").ClearContents()
# Row 122
$ws.Range("H122").Value = 1998.5
$ws.Range("I122").Value = 1998.5
$ws.Range("K122").Value = 5995.5
$ws.Range("M122").Value = -3545.5
# Row 132
$ws.Range("H132").Value = 2571.7368
$ws.Range("I132").Value = 2439.2856
$ws.Range("K132").Value = 7317.8568
$ws.Range("M132").Value = -4787.8568
# Row 136
$ws.Range("H136").Value = 5727
$ws.Range("I136").Value = 3699.6667
$ws.Range("J136").Value = 7754.3335
$ws.Range("K136").Value = 11099.0001
$ws.Range("L136").Value = 23263.0005
$ws.Range("M136").Value = -8549.000100000001
$ws.Range("N136").Value = -28363.0005

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2102.6072
$ws.Range("I105").Value = 2054.96
$ws.Range("K105").Value = 2054.96
$ws.Range("M105").Value = -307.96
# Row 107
$ws.Range("H107").Value = 1533.625
$ws.Range("I107").Value = 1424
$ws.Range("J107").Value = 1799.8572
$ws.Range("K107").Value = 1424
$ws.Range("L107").Value = 1799.8572
$ws.Range("M107").Value = 496
$ws.Range("N107").Value = -5639.8572

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1976.579
$ws.Range("I31").Value = 1327.6666
$ws.Range("K31").Value = 1327.6666
$ws.Range("M31").Value = -1032.6666
# Row 34
$ws.Range("H34").Value = 1976.579
$ws.Range("I34").Value = 1327.6666
$ws.Range("K34").Value = 1327.6666
$ws.Range("M34").Value = -1125.6666
# Row 86
$ws.Range("H86").Value = 2121.111
$ws.Range("I86").Value = 1665.3334
$ws.Range("K86").Value = 1665.3334
$ws.Range("M86").Value = -542.3334
# Row 89
$ws.Range("H89").Value = 2121.111
$ws.Range("I89").Value = 1665.3334
$ws.Range("K89").Value = 8326.666999999999
$ws.Range("M89").Value = -2710.666999999999
# Row 132
$ws.Range("H132").Value = 2930
$ws.Range("I132").Value = 1349.875
$ws.Range("J132").Value = 4334.5557
$ws.Range("K132").Value = 4049.625
$ws.Range("L132").Value = 13003.6671
$ws.Range("M132").Value = -1519.625
$ws.Range("N132").Value = -18063.6671

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 43
$ws.Range("I8").Value = 43
$ws.Range("K8").Value = 129
$ws.Range("M8").Value = 10
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
# Row 127
$ws.Range("H127").Value = 1900
$ws.Range("J127").Value = 1900
$ws.Range("L127").Value = 5700
$ws.Range("N127").Value = -15620
# Row 131
$ws.Range("H131").Value = 785.8081
$ws.Range("J131").Value = 794.73956
$ws.Range("L131").Value = 2384.21868
$ws.Range("N131").Value = -12464.21868
# Row 137
$ws.Range("H137").Value = 3456.0527
$ws.Range("I137").Value = 1666.6666
$ws.Range("J137").Value = 3791.5625
$ws.Range("K137").Value = 4999.9998
$ws.Range("L137").Value = 11374.6875
$ws.Range("M137").Value = 100.0002000000004
$ws.Range("N137").Value = -21574.6875
# Row 140
$ws.Range("H140").Value = 1485.5588
$ws.Range("I140").Value = 823.65
$ws.Range("J140").Value = 2431.1428
$ws.Range("K140").Value = 2470.95
$ws.Range("L140").Value = 7293.428400000001
$ws.Range("M140").Value = 2709.05
$ws.Range("N140").Value = -17653.4284
# Row 141
$ws.Range("H141").Value = 3359.7144
$ws.Range("I141").Value = 3359.7144
$ws.Range("K141").Value = 10079.1432
$ws.Range("M141").Value = -4899.143199999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3720.125
$ws.Range("I80").Value = 3552.8
$ws.Range("K80").Value = 3552.8
$ws.Range("M80").Value = -2554.8
# Row 83
$ws.Range("H83").Value = 3720.125
$ws.Range("I83").Value = 3552.8
$ws.Range("K83").Value = 17764
$ws.Range("M83").Value = -12772
# Row 107
$ws.Range("H107").Value = 867.6667
$ws.Range("J107").Value = 2103
$ws.Range("L107").Value = 2103
$ws.Range("N107").Value = -5943
# Row 113
$ws.Range("H113").Value = 1299.75
$ws.Range("I113").Value = 1099.5
$ws.Range("K113").Value = 1099.5
$ws.Range("M113").Value = 1070.5
# Row 122
$ws.Range("H122").Value = 1669.1578
$ws.Range("I122").Value = 1259.75
$ws.Range("K122").Value = 3779.25
$ws.Range("M122").Value = -1329.25
# Row 127
$ws.Range("H127").Value = 27704.666
$ws.Range("J127").Value = 27704.666
$ws.Range("L127").Value = 27704.666
$ws.Range("N127").Value = -37624.666
# Row 132
$ws.Range("H132").Value = 2026747.8
$ws.Range("I132").Value = 4275283.5
$ws.Range("J132").Value = 3065.8
$ws.Range("K132").Value = 12825850.5
$ws.Range("L132").Value = 9197.400000000001
$ws.Range("M132").Value = -12823320.5
$ws.Range("N132").Value = -14257.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1668.2941
$ws.Range("I22").Value = 1545.8889
$ws.Range("J22").Value = 1806
$ws.Range("K22").Value = 1545.8889
$ws.Range("L22").Value = 1806
$ws.Range("M22").Value = -1250.8889
$ws.Range("N22").Value = -2396
# Row 27
$ws.Range("H27").Value = 1668.2941
$ws.Range("I27").Value = 1545.8889
$ws.Range("J27").Value = 1806
$ws.Range("K27").Value = 1545.8889
$ws.Range("L27").Value = 1806
$ws.Range("M27").Value = -1438.8889
$ws.Range("N27").Value = -2020
# Row 122
$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 136
$ws.Range("H136").Value = 2413.5557
$ws.Range("I136").Value = 1433.3
$ws.Range("K136").Value = 4299.9
$ws.Range("M136").Value = -1749.9

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 61799.2
$ws.Range("J108").Value = 61799.2
$ws.Range("L108").Value = 61799.2
$ws.Range("N108").Value = -69479.2
# Row 122
$ws.Range("H122").Value = 66499.914
$ws.Range("I122").Value = 112253.57
$ws.Range("K122").Value = 336760.71
$ws.Range("M122").Value = -334310.71
